# Edit script: split the intro sentence into several runs describing the
# additional "waveArithUnit.do" include, and add a new paragraph about
# TimingArithUnit.do. Also bring the styles.xml latent-style exception
# list up to date (Normal Table / Table Web 2 / Table Theme entries).

$d = $word.ActiveDocument

$apos = [char]0x2019

# ---------------------------------------------------------------------
# 1) Rework the first paragraph's text run into five runs, add a blank
#    paragraph, and add a new paragraph with the TimingArithUnit note.
# ---------------------------------------------------------------------
$oldPara = '<w:r><w:t>Include LogicUnitWave.do! I' + $apos + 've made changes</w:t></w:r>'

$newPara = '<w:r><w:t xml:space="preserve">Include </w:t></w:r>' +
           '<w:r><w:t>waveLogicUnit</w:t></w:r>' +
           '<w:r><w:t>.do</w:t></w:r>' +
           '<w:r><w:t xml:space="preserve"> and waveArithUnit.do</w:t></w:r>' +
           '<w:r><w:t>! I' + $apos + 've made changes</w:t></w:r>' +
           '</w:p><w:p/><w:p><w:r><w:t>Changed TimingArithUnit.do to run for 15600 ns.</w:t></w:r>'

$wxml = $d.WordOpenXML

if ($wxml.IndexOf($oldPara) -lt 0) {
    throw "Could not locate the target paragraph text in WordOpenXML"
}
$wxml = $wxml.Replace($oldPara, $newPara)

# ---------------------------------------------------------------------
# 2) Add the three missing lsdException entries to styles.xml
# ---------------------------------------------------------------------
$htmlVariable = '<w:lsdException w:name="HTML Variable" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$normalTable  = '<w:lsdException w:name="Normal Table" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$wxml = $wxml.Replace($htmlVariable, $htmlVariable + $normalTable)

$tableWeb1 = '<w:lsdException w:name="Table Web 1" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$tableWeb2 = '<w:lsdException w:name="Table Web 2" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$wxml = $wxml.Replace($tableWeb1, $tableWeb1 + $tableWeb2)

$tableGrid  = '<w:lsdException w:name="Table Grid" w:uiPriority="39"/>'
$tableTheme = '<w:lsdException w:name="Table Theme" w:semiHidden="1" w:unhideWhenUsed="1"/>'
$wxml = $wxml.Replace($tableGrid, $tableGrid + $tableTheme)

# ---------------------------------------------------------------------
# 3) Push the rewritten package back into the document
# ---------------------------------------------------------------------
$d.Range().InsertXML($wxml)

Write-Output "done"
